$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Kevin Roycroft <Kevin.Roycroft@ftsi.com>; Saul Martinez <Saul.Martinez@ftsi.com>; Issac Jacquez <Issac.Jacquez@ftsi.com>; Andre Malone <Andre.Malone@ftsi.com>; Charles Hale <Charles.Hale@ftsi.com>; James Demoss <James.Demoss@ftsi.com>; Jose Santos <Jose.Santos@ftsi.com>; Ziad Kanaan <Ziad.Kanaan@ftsi.com>'
$ws.Range("A2").Value = 'Max Taylor <Max.Taylor@ftsi.com>; Teddy Woodson <Teddy.Woodson@ftsi.com>; Daniel Jones <Daniel.Jones3@ftsi.com>; Willie Atkinson <Willie.Atkinson@ftsi.com>; Odessa Technical Advisors <OdessaTechAdvisors@ftsi.com>; ''tbritton@kcftech.com''; Richard Bushnell <Richard.Bushnell@ftsi.com>'
$ws.Range("A3").Value = 'Dedron Sells <Dedron.Sells@ftsi.com>; Joseph Esposito <Joseph.Esposito@ftsi.com>; Joshua Haynes <Joshua.Haynes@ftsi.com>; John Romiluyi <John.Romiluyi@ftsi.com>; Sammy Campos <Samuel.Campos@ftsi.com>; Mohammed Barakat <Mohammed.Barakat@ftsi.com>'
$ws.Range("A4").Value = 'Max Taylor <Max.Taylor@ftsi.com>; Teddy Woodson <Teddy.Woodson@ftsi.com>; Daniel Jones <Daniel.Jones3@ftsi.com>; Willie Atkinson <Willie.Atkinson@ftsi.com>; Odessa Technical Advisors <OdessaTechAdvisors@ftsi.com>; tbritton@kcftech.com; Richard Bushnell <Richard.Bushnell@ftsi.com>'
$ws.Range("A5").Value = 'Raul Galvan-Kalinchuk <Raul.Galvan-Kalinchuk@ftsi.com>; Alberto Ortiz <Alberto.Ortiz@ftsi.com>; Gustavo Machuca <Gustavo.Machuca@ftsi.com>; Jose Quintanilla <Jose.Quintanilla@ftsi.com>; Miguel Zuniga <Miguel.Zuniga@ftsi.com>; Marco Carmona <Marco.Carmona@ftsi.com>; Brandon Lowrey <Brandon.Lowrey@ftsi.com>; Dave Perkins <Dave.Perkins@ftsi.com>; Jose Martinez <Jose.Martinez2@ftsi.com>'
$ws.Range("A6").Value = 'Max Taylor <Max.Taylor@ftsi.com>; Teddy Woodson <Teddy.Woodson@ftsi.com>; Daniel Jones <Daniel.Jones3@ftsi.com>; Willie Atkinson <Willie.Atkinson@ftsi.com>; Odessa Technical Advisors <OdessaTechAdvisors@ftsi.com>; ''tbritton@kcftech.com''; Richard Bushnell <Richard.Bushnell@ftsi.com>'
$ws.Range("A7").Value = 'Byron Austin <Byron.Austin@ftsi.com>; Hussein Aydi <Hussein.Aydi@ftsi.com>; Jesus Castro <Jesus.Castro@ftsi.com>; Juan Olazaba <Juan.Olazaba@ftsi.com>; Brandon Slaton <Brandon.Slaton@ftsi.com>; Nicholas Perez-Lasala <Nicholas.Perez@ftsi.com>; Michael Driver <Michael.Driver@ftsi.com>; Aaron Getgen <Aaron.Getgen@ftsi.com>; John Demoss <John.Demoss@ftsi.com>'
$ws.Range("A8").Value = 'Max Taylor <Max.Taylor@ftsi.com>; Teddy Woodson <Teddy.Woodson@ftsi.com>; Daniel Jones <Daniel.Jones3@ftsi.com>; Willie Atkinson <Willie.Atkinson@ftsi.com>; Odessa Technical Advisors <OdessaTechAdvisors@ftsi.com>; ''tbritton@kcftech.com''; Richard Bushnell <Richard.Bushnell@ftsi.com>'
$ws.Range("A9").Value = 'Jose Rodriguez <Jose.Rodriguez5@ftsi.com>; Jay Crawford <Jay.Crawford@ftsi.com>; Clark Akers <Clark.Akers@ftsi.com>; Saul Lopez <Saul.Lopez@ftsi.com>; Gerardo Salinas <Gerardo.Salinas@ftsi.com>; Angel Alvarado <Angel.Alvarado@ftsi.com>, Jose Delgado <Jose.Delgado2@ftsi.com>, Max Taylor <Max.Taylor@ftsi.com>;William Blunk <William.Blunk@ftsi.com>'
$ws.Range("A10").Value = 'Cody Irwin <Cody.Irwin@ftsi.com>; Flavio Bombonato <fbombonato@kcftech.com>; KCF Tech <sentryftsi@kcftech.com>; Jose Quintanilla <Jose.Quintanilla@ftsi.com>; Jonathan Woodson <Jonathan.Woodson@ftsi.com>;Richard Bushnell <Richard.Bushnell@ftsi.com>;Willie Atkinson <Willie.Atkinson@ftsi.com>; Daniel Jones <Daniel.Jones3@ftsi.com>'
$ws.Range("A11").Value = 'Luis Salinas <Luis.Salinas@ftsi.com>; Fernando Segovia <Fernando.Segovia@ftsi.com>; Roberto Rivera <Roberto.Rivera@ftsi.com>; Jesus Barrera <Jesus.Barrera2@ftsi.com>; Ismael Soto <Ismael.Soto@ftsi.com>; Jeremy Grissom <Jeremy.Grissom@ftsi.com>; Ryan Castillo <Ryan.Castillo@ftsi.com>; Donovan Denmon <Donovan.Denmon@ftsi.com>'
$ws.Range("A12").Value = 'Max Taylor <Max.Taylor@ftsi.com>; Teddy Woodson <Teddy.Woodson@ftsi.com>; Daniel Jones <Daniel.Jones3@ftsi.com>; Willie Atkinson <Willie.Atkinson@ftsi.com>; Odessa Technical Advisors <OdessaTechAdvisors@ftsi.com>; ''tbritton@kcftech.com''; Richard Bushnell <Richard.Bushnell@ftsi.com>'

$ws.Range("C21").Select()
